$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 73-74; this pushes the existing rows 73..180
# down to 75..182 (matching the target dimension A1:T182).
$ws.Rows("73:74").Insert()

# Row 73 - new record (Primera, Cultivar IV Region, 2024 season)
$ws.Range("A73").Value = 5
$ws.Range("B73").Value = "Macroferia Regional de Talca"
$ws.Range("C73").Value = "Maule"
$ws.Range("D73").Value = 45219
$ws.Range("E73").Value = 7
$ws.Range("F73").Value = "Fruta"
$ws.Range("G73").Value = 100107
$ws.Range("H73").Value = "Otros"
$ws.Range("I73").Value = 100107002
$ws.Range("J73").Value = "Chirimoya"
$ws.Range("K73").Value = "Cultivar IV Región"
$ws.Range("L73").Value = "Primera"
$ws.Range("M73").Value = 230
$ws.Range("N73").Value = 20000
$ws.Range("O73").Value = 20000
$ws.Range("P73").Value = 20000
$ws.Range("Q73").Value = "`$/bandeja 10 kilos"
$ws.Range("R73").Value = "Provincia de Limarí"
$ws.Range("S73").Value = 2000
$ws.Range("T73").Value = 10

# Row 74 - new record (Segunda, Cultivar IV Region, 2024 season)
$ws.Range("A74").Value = 5
$ws.Range("B74").Value = "Macroferia Regional de Talca"
$ws.Range("C74").Value = "Maule"
$ws.Range("D74").Value = 45219
$ws.Range("E74").Value = 7
$ws.Range("F74").Value = "Fruta"
$ws.Range("G74").Value = 100107
$ws.Range("H74").Value = "Otros"
$ws.Range("I74").Value = 100107002
$ws.Range("J74").Value = "Chirimoya"
$ws.Range("K74").Value = "Cultivar IV Región"
$ws.Range("L74").Value = "Segunda"
$ws.Range("M74").Value = 200
$ws.Range("N74").Value = 18000
$ws.Range("O74").Value = 18000
$ws.Range("P74").Value = 18000
$ws.Range("Q74").Value = "`$/bandeja 10 kilos"
$ws.Range("R74").Value = "Provincia de Limarí"
$ws.Range("S74").Value = 1800
$ws.Range("T74").Value = 10
